# Update the weekly Pomelo (Lo Valledor) price sheet.
# The underlying data rows were refreshed for the new reporting week:
#  - existing rows 2-20 were updated in place (date, quality, volume,
#    min/max/weighted price, unit, origin region and $/Kg as needed)
#  - a brand-new observation was appended as row 21
# Values below were derived from the published row-by-row diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Date number format used by column D (reuse existing style format code)
$dateFmt = $ws.Cells.Item(2, 4).NumberFormat

# Row 2
$ws.Cells.Item(2, 4).Value = 44356
$ws.Cells.Item(2, 4).NumberFormat = $dateFmt
$ws.Cells.Item(2, 13).Value = 24
$ws.Cells.Item(2, 14).Value = 200000
$ws.Cells.Item(2, 15).Value = 230000
$ws.Cells.Item(2, 16).Value = 215000
$ws.Cells.Item(2, 18).Value = 'Región Metropolitana'
$ws.Cells.Item(2, 19).Value = 614

# Row 3
$ws.Cells.Item(3, 4).Value = 44208
$ws.Cells.Item(3, 4).NumberFormat = $dateFmt
$ws.Cells.Item(3, 12).Value = 'Primera'
$ws.Cells.Item(3, 13).Value = 16
$ws.Cells.Item(3, 14).Value = 180000
$ws.Cells.Item(3, 15).Value = 180000
$ws.Cells.Item(3, 16).Value = 180000
$ws.Cells.Item(3, 19).Value = 514

# Row 4
$ws.Cells.Item(4, 4).Value = 44389
$ws.Cells.Item(4, 4).NumberFormat = $dateFmt
$ws.Cells.Item(4, 12).Value = 'Especial'
$ws.Cells.Item(4, 13).Value = 18
$ws.Cells.Item(4, 14).Value = 200000
$ws.Cells.Item(4, 15).Value = 200000
$ws.Cells.Item(4, 16).Value = 200000
$ws.Cells.Item(4, 17).Value = '$/bins (350 kilos)'
$ws.Cells.Item(4, 18).Value = 'Provincia de Quillota'
$ws.Cells.Item(4, 19).Value = 571
$ws.Cells.Item(4, 20).Value = 350

# Row 5
$ws.Cells.Item(5, 4).Value = 44193
$ws.Cells.Item(5, 4).NumberFormat = $dateFmt
$ws.Cells.Item(5, 13).Value = 8
$ws.Cells.Item(5, 14).Value = 150000
$ws.Cells.Item(5, 15).Value = 150000
$ws.Cells.Item(5, 16).Value = 150000
$ws.Cells.Item(5, 19).Value = 429

# Row 6
$ws.Cells.Item(6, 4).Value = 44298
$ws.Cells.Item(6, 4).NumberFormat = $dateFmt
$ws.Cells.Item(6, 12).Value = 'Especial'
$ws.Cells.Item(6, 13).Value = 15
$ws.Cells.Item(6, 14).Value = 450000
$ws.Cells.Item(6, 15).Value = 450000
$ws.Cells.Item(6, 16).Value = 450000
$ws.Cells.Item(6, 19).Value = 1286

# Row 7
$ws.Cells.Item(7, 4).Value = 44298
$ws.Cells.Item(7, 4).NumberFormat = $dateFmt
$ws.Cells.Item(7, 14).Value = 430000
$ws.Cells.Item(7, 15).Value = 430000
$ws.Cells.Item(7, 16).Value = 430000
$ws.Cells.Item(7, 18).Value = 'Región Metropolitana'
$ws.Cells.Item(7, 19).Value = 1229

# Row 8
$ws.Cells.Item(8, 4).Value = 44201
$ws.Cells.Item(8, 4).NumberFormat = $dateFmt
$ws.Cells.Item(8, 12).Value = 'Especial'
$ws.Cells.Item(8, 13).Value = 8
$ws.Cells.Item(8, 14).Value = 200000
$ws.Cells.Item(8, 15).Value = 200000
$ws.Cells.Item(8, 16).Value = 200000
$ws.Cells.Item(8, 17).Value = '$/bins (350 kilos)'
$ws.Cells.Item(8, 19).Value = 571
$ws.Cells.Item(8, 20).Value = 350

# Row 9
$ws.Cells.Item(9, 4).Value = 44201
$ws.Cells.Item(9, 4).NumberFormat = $dateFmt
$ws.Cells.Item(9, 12).Value = 'Primera'
$ws.Cells.Item(9, 13).Value = 16
$ws.Cells.Item(9, 14).Value = 170000
$ws.Cells.Item(9, 15).Value = 170000
$ws.Cells.Item(9, 16).Value = 170000
$ws.Cells.Item(9, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(9, 19).Value = 486

# Row 10
$ws.Cells.Item(10, 4).Value = 44308
$ws.Cells.Item(10, 4).NumberFormat = $dateFmt
$ws.Cells.Item(10, 14).Value = 280000
$ws.Cells.Item(10, 15).Value = 280000
$ws.Cells.Item(10, 16).Value = 280000
$ws.Cells.Item(10, 19).Value = 800

# Row 11
$ws.Cells.Item(11, 4).Value = 44363
$ws.Cells.Item(11, 4).NumberFormat = $dateFmt
$ws.Cells.Item(11, 13).Value = 20
$ws.Cells.Item(11, 14).Value = 200000
$ws.Cells.Item(11, 15).Value = 230000
$ws.Cells.Item(11, 16).Value = 215000
$ws.Cells.Item(11, 18).Value = 'Provincia de Limarí'
$ws.Cells.Item(11, 19).Value = 614

# Row 12
$ws.Cells.Item(12, 4).Value = 44189
$ws.Cells.Item(12, 4).NumberFormat = $dateFmt
$ws.Cells.Item(12, 11).Value = 'Start Ruby'
$ws.Cells.Item(12, 13).Value = 16
$ws.Cells.Item(12, 14).Value = 150000
$ws.Cells.Item(12, 15).Value = 150000
$ws.Cells.Item(12, 16).Value = 150000
$ws.Cells.Item(12, 19).Value = 429

# Row 13
$ws.Cells.Item(13, 4).Value = 44312
$ws.Cells.Item(13, 4).NumberFormat = $dateFmt
$ws.Cells.Item(13, 12).Value = 'Segunda'
$ws.Cells.Item(13, 13).Value = 10
$ws.Cells.Item(13, 14).Value = 330000
$ws.Cells.Item(13, 15).Value = 330000
$ws.Cells.Item(13, 16).Value = 330000
$ws.Cells.Item(13, 18).Value = 'Región Metropolitana'
$ws.Cells.Item(13, 19).Value = 943

# Row 15
$ws.Cells.Item(15, 4).Value = 44309
$ws.Cells.Item(15, 4).NumberFormat = $dateFmt
$ws.Cells.Item(15, 12).Value = 'Primera'
$ws.Cells.Item(15, 13).Value = 16
$ws.Cells.Item(15, 14).Value = 350000
$ws.Cells.Item(15, 15).Value = 350000
$ws.Cells.Item(15, 16).Value = 350000
$ws.Cells.Item(15, 18).Value = 'Región Metropolitana'
$ws.Cells.Item(15, 19).Value = 1000

# Row 16
$ws.Cells.Item(16, 4).Value = 44196
$ws.Cells.Item(16, 4).NumberFormat = $dateFmt
$ws.Cells.Item(16, 11).Value = 'Red Blush'
$ws.Cells.Item(16, 13).Value = 12
$ws.Cells.Item(16, 14).Value = 130000
$ws.Cells.Item(16, 15).Value = 130000
$ws.Cells.Item(16, 16).Value = 130000
$ws.Cells.Item(16, 18).Value = 'Provincia de Limarí'
$ws.Cells.Item(16, 19).Value = 371

# Row 17
$ws.Cells.Item(17, 4).Value = 44446
$ws.Cells.Item(17, 4).NumberFormat = $dateFmt
$ws.Cells.Item(17, 13).Value = 14
$ws.Cells.Item(17, 14).Value = 150000
$ws.Cells.Item(17, 15).Value = 160000
$ws.Cells.Item(17, 16).Value = 155000
$ws.Cells.Item(17, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(17, 19).Value = 443

# Row 20
$ws.Cells.Item(20, 4).Value = 44167
$ws.Cells.Item(20, 4).NumberFormat = $dateFmt
$ws.Cells.Item(20, 13).Value = 140
$ws.Cells.Item(20, 14).Value = 9800
$ws.Cells.Item(20, 15).Value = 9800
$ws.Cells.Item(20, 16).Value = 9800
$ws.Cells.Item(20, 17).Value = '$/caja 14 kilos empedrada'
$ws.Cells.Item(20, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(20, 19).Value = 700
$ws.Cells.Item(20, 20).Value = 14

# Row 21
$ws.Cells.Item(21, 1).Value = 6
$ws.Cells.Item(21, 2).Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Cells.Item(21, 3).Value = 'Metropolitana'
$ws.Cells.Item(21, 4).Value = 44400
$ws.Cells.Item(21, 4).NumberFormat = $dateFmt
$ws.Cells.Item(21, 5).Value = 13
$ws.Cells.Item(21, 6).Value = 'Fruta'
$ws.Cells.Item(21, 7).Value = 100102
$ws.Cells.Item(21, 8).Value = 'Cítricos'
$ws.Cells.Item(21, 9).Value = 100102006
$ws.Cells.Item(21, 10).Value = 'Pomelo'
$ws.Cells.Item(21, 11).Value = 'Start Ruby'
$ws.Cells.Item(21, 12).Value = 'Primera'
$ws.Cells.Item(21, 13).Value = 140
$ws.Cells.Item(21, 14).Value = 9800
$ws.Cells.Item(21, 15).Value = 9800
$ws.Cells.Item(21, 16).Value = 9800
$ws.Cells.Item(21, 17).Value = '$/caja 14 kilos empedrada'
$ws.Cells.Item(21, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(21, 19).Value = 700
$ws.Cells.Item(21, 20).Value = 14

